# Update Leave Card 12/22/2023 10:59 AM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Row 57: period end-date corrected from 1/1/2023 to 1/31/2023
$ws.Range("A57").Value = 44957

# Rows 58-66: fill in the monthly period dates and the 1.25 VL earned amount
$ws.Range("A58").Value = 44985
$ws.Range("C58").Value = 1.25

$ws.Range("A59").Value = 45016
$ws.Range("C59").Value = 1.25

$ws.Range("A60").Value = 45046
$ws.Range("C60").Value = 1.25

$ws.Range("A61").Value = 45077
$ws.Range("C61").Value = 1.25

$ws.Range("A62").Value = 45107
$ws.Range("C62").Value = 1.25

$ws.Range("A63").Value = 45138
$ws.Range("C63").Value = 1.25

$ws.Range("A64").Value = 45169
$ws.Range("C64").Value = 1.25

$ws.Range("A65").Value = 45199
$ws.Range("C65").Value = 1.25

$ws.Range("A66").Value = 45230
$ws.Range("C66").Value = 1.25

# Row 67: leave taken entry (Paternity Leave) with remarks date range
$ws.Range("A67").Value = 45260
$ws.Range("B67").Value = "PL(7-0-0)"
$ws.Range("C67").Value = 1.25
$ws.Range("K67").Value = "11/23,24,25,26,30, 12/1,2/2023"

# Rows 68-71: continue the monthly period dates
$ws.Range("A68").Value = 45291
$ws.Range("A69").Value = 45322
$ws.Range("A70").Value = 45351
$ws.Range("A71").Value = 45382
